$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column C
$ws.Range("C1").Value = "Row"

# Fill column C values: rows 2-15 -> 1, rows 16-19 -> 2
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 3).Value = 1
}
for ($r = 16; $r -le 19; $r++) {
    $ws.Cells.Item($r, 3).Value = 2
}

# Update the selection to match the new target cell
$ws.Range("C20").Select()
